$d = $word.ActiveDocument

# The document contains a single table listing "Đề tài / Học viên / Giảng viên".
# Row 1 is the header, row 2 and row 3 are data rows.
$t = $d.Tables.Item(1)

# Update the first data row (row 2) in place.
$t.Cell(2, 1).Range.Text = "de tai"
$t.Cell(2, 2).Range.Text = "hv01"
$t.Cell(2, 3).Range.Text = "giảng viên 2"

# Remove the second data row (row 3) entirely.
$t.Rows.Item(3).Delete()
